$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# 1. "PostgresSQL" textbox (object 19) -> "MySQL": only the first run's
#    text ("Postgres") is replaced with "My"; the second run ("SQL") is
#    left untouched so the two runs / formatting stay intact.
$sh = $s.Shapes.Item(9)
$tf = $sh.TextFrame
$tr = $tf.TextRange

$origHeight = $sh.Height

$run = $tr.Characters(1, 8)
$run.Text = "My"

# The shape auto-fits its height to the text (spAutoFit); editing the
# text recalculates that height, so restore the original value. (The
# shape model stores Height in points at single precision, so nudge
# slightly off the exact point value to land back on the same EMU.)
$sh.Height = $origHeight + 0.00001

# 2. Remove the plain white rectangle (object 21) that sat on top of the
#    table in the top-right corner of the slide.
$s.Shapes.Item("object 21").Delete()
